# Recolor the highlighted terms on slide 2 (first body placeholder):
#   - the three "yellow" (FFFF00) highlighted runs become a softer
#     yellow (FFE599): "conductive shielding cap ", "configured",
#     " to be assembled to"
#   - the "orange" (FF9900) highlighted run on "case body" becomes a
#     softer orange (F9CB9C)
# The trailing closing-quote run and the paragraph's endParaRPr (also
# FF9900) are left untouched, matching the source diff.

function ToRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$softYellow = ToRGB 0xFF 0xE5 0x99   # FFE599
$softOrange = ToRGB 0xF9 0xCB 0x9C   # F9CB9C

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Run boundaries (1-based char index, length) inside the paragraph:
#   1: "The case being deprived of a shielding cap, “"   (1,45)  - no highlight
#   2: "said connector comprises a "                      (46,27) - no highlight
#   3: "conductive shielding cap "                        (73,25) - FFFF00 -> FFE599
#   4: "configured"                                       (98,10) - FFFF00 -> FFE599
#   5: " to be assembled to"                               (108,19) - FFFF00 -> FFE599
#   6: " the "                                             (127,5)  - no highlight
#   7: "case body"                                         (132,9)  - FF9900 -> F9CB9C
#   8: "”"                                                 (141,1)  - no highlight

$tr.Characters(73, 25).Font.Highlight.RGB = $softYellow
$tr.Characters(98, 10).Font.Highlight.RGB = $softYellow
$tr.Characters(108, 19).Font.Highlight.RGB = $softYellow
$tr.Characters(132, 9).Font.Highlight.RGB = $softOrange
